$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells keep their original "Text" storage type (the source file stores
# these look-like-numbers as inline strings) instead of Excel auto-coercing the
# assigned strings into Number/Percentage values.
$cells = @{
    "D2" = "328.94"
    "E2" = "1.59%"
    "D3" = "41.26"
    "E3" = "5.13%"
    "D4" = "5.648"
    "E4" = "-0.68%"
    "D5" = "0.08188"
    "E5" = "2.27%"
    "D6" = "8.742"
    "E6" = "1.34%"
    "D7" = "2.006"
    "E7" = "0.24%"
    "E8" = "-0.18%"
    "E9" = "1.45%"
    "D10" = "0.9200"
    "E10" = "-0.76%"
    "D11" = "0.1276"
    "E11" = "2.93%"
    "D12" = "0.1948"
    "E12" = "-1.76%"
    "D13" = "0.09228"
    "E13" = "0.22%"
    "D14" = "0.03886"
    "E14" = "7.71%"
    "E15" = "0.89%"
    "D16" = "0.001297"
    "E16" = "0.43%"
    "D17" = "0.006191"
    "E17" = "1.00%"
    "D19" = "3.447"
    "E19" = "2.93%"
    "E20" = "0.34%"
    "D21" = "8.221"
    "E21" = "-5.38%"
    "E22" = "0.04%"
    "D23" = "0.2411"
    "E23" = "-0.07%"
    "D24" = "0.04393"
    "E24" = "-0.43%"
    "E25" = "-0.59%"
    "D26" = "0.004308"
    "E26" = "-6.87%"
    "E27" = "4.38%"
    "D39" = "0.02794"
    "E39" = "12.07%"
    "D40" = "0.05404"
    "E40" = "0.96%"
    "D41" = "0.007698"
    "E41" = "3.11%"
    "D42" = "0.1417"
    "E42" = "0.86%"
    "D43" = "0.008947"
    "E43" = "-6.97%"
    "D44" = "0.002172"
    "E44" = "2.59%"
    "D45" = "0.01157"
    "E45" = "12.23%"
    "D46" = "0.00006585"
    "E46" = "-2.00%"
    "E47" = "0.03%"
    "D48" = "0.003205"
    "E48" = "7.83%"
    "D49" = "0.002281"
    "E49" = "-0.46%"
    "E50" = "0.03%"
    "E51" = "0.03%"
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}
